$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style/format of the last existing row (73) down through the new rows (74-103)
$ws.Range("A73:F73").Copy()
$ws.Range("A74:F103").PasteSpecial(-4122)

$ws.Cells.Item(74, 1).Value = 45222.35347267361
$ws.Cells.Item(74, 2).Value = "Raden Dewi Ghita Ghaida P"
$ws.Cells.Item(74, 3).Value = 23123010
$ws.Cells.Item(74, 4).Value = "Regular"
$ws.Cells.Item(74, 5).Value = 160
$ws.Cells.Item(74, 6).Value = 38

$ws.Cells.Item(75, 1).Value = 45222.35351905093
$ws.Cells.Item(75, 2).Value = "Sunita Hiu Lienshin"
$ws.Cells.Item(75, 3).Value = 13120132
$ws.Cells.Item(75, 4).Value = "Regular"
$ws.Cells.Item(75, 5).Value = 165
$ws.Cells.Item(75, 6).Value = 29

$ws.Cells.Item(76, 1).Value = 45222.362550439815
$ws.Cells.Item(76, 2).Value = "Wahyu Eko Prasetyo Akbar"
$ws.Cells.Item(76, 3).Value = 23123043
$ws.Cells.Item(76, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(76, 5).Value = 172
$ws.Cells.Item(76, 6).Value = 34

$ws.Cells.Item(77, 1).Value = 45222.371023125
$ws.Cells.Item(77, 2).Value = "Pradana Jayawardana"
$ws.Cells.Item(77, 3).Value = 23123025
$ws.Cells.Item(77, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(77, 5).Value = 177
$ws.Cells.Item(77, 6).Value = 71

$ws.Cells.Item(78, 1).Value = 45222.371432175925
$ws.Cells.Item(78, 2).Value = "DIAN AKBAR KARISMASANI"
$ws.Cells.Item(78, 3).Value = 23123031
$ws.Cells.Item(78, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(78, 5).Value = 174
$ws.Cells.Item(78, 6).Value = 41

$ws.Cells.Item(79, 1).Value = 45222.38621003472
$ws.Cells.Item(79, 2).Value = "Baradiant Ivano Leotman"
$ws.Cells.Item(79, 3).Value = 23123032
$ws.Cells.Item(79, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(79, 5).Value = 175
$ws.Cells.Item(79, 6).Value = 40

$ws.Cells.Item(80, 1).Value = 45222.39446164352
$ws.Cells.Item(80, 2).Value = "ramadani putra"
$ws.Cells.Item(80, 3).Value = 23123016
$ws.Cells.Item(80, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(80, 5).Value = 160
$ws.Cells.Item(80, 6).Value = 31

$ws.Cells.Item(81, 1).Value = 45222.415436655094
$ws.Cells.Item(81, 2).Value = "Rizky Andri Nurachman"
$ws.Cells.Item(81, 3).Value = 23123020
$ws.Cells.Item(81, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(81, 5).Value = 172
$ws.Cells.Item(81, 6).Value = 59

$ws.Cells.Item(82, 1).Value = 45222.432648807866
$ws.Cells.Item(82, 2).Value = "Henry Situmorang "
$ws.Cells.Item(82, 3).Value = 23123034
$ws.Cells.Item(82, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(82, 5).Value = 165
$ws.Cells.Item(82, 6).Value = "32 WPM (150 CPM)"

$ws.Cells.Item(83, 1).Value = 45222.44384939814
$ws.Cells.Item(83, 2).Value = "M. RAMADHAN SAPUTRA"
$ws.Cells.Item(83, 3).Value = 23123038
$ws.Cells.Item(83, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(83, 5).Value = 178
$ws.Cells.Item(83, 6).Value = 43

$ws.Cells.Item(84, 1).Value = 45222.454839872684
$ws.Cells.Item(84, 2).Value = "riza afandi"
$ws.Cells.Item(84, 3).Value = 20023009
$ws.Cells.Item(84, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(84, 5).Value = 171
$ws.Cells.Item(84, 6).Value = 34

$ws.Cells.Item(85, 1).Value = 45222.47543734954
$ws.Cells.Item(85, 2).Value = "Ni Luh Putu Andrea Maurilla Sarasvanya"
$ws.Cells.Item(85, 3).Value = 23123027
$ws.Cells.Item(85, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(85, 5).Value = 155
$ws.Cells.Item(85, 6).Value = 47

$ws.Cells.Item(86, 1).Value = 45222.47562165509
$ws.Cells.Item(86, 2).Value = "Meidiono Untoro"
$ws.Cells.Item(86, 3).Value = 23123049
$ws.Cells.Item(86, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(86, 5).Value = 170
$ws.Cells.Item(86, 6).Value = 46

$ws.Cells.Item(87, 1).Value = 45222.477326249995
$ws.Cells.Item(87, 2).Value = "RIDO PAHOTAN TUA MANIK"
$ws.Cells.Item(87, 3).Value = 23123022
$ws.Cells.Item(87, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(87, 5).Value = 166
$ws.Cells.Item(87, 6).Value = 32

$ws.Cells.Item(88, 1).Value = 45222.55868618056
$ws.Cells.Item(88, 2).Value = "Muhammad Siddiq B"
$ws.Cells.Item(88, 3).Value = 23123021
$ws.Cells.Item(88, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(88, 5).Value = 180
$ws.Cells.Item(88, 6).Value = 42

$ws.Cells.Item(89, 1).Value = 45222.59568787037
$ws.Cells.Item(89, 2).Value = "Tiyas Sinta Rahmania"
$ws.Cells.Item(89, 3).Value = 23123026
$ws.Cells.Item(89, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(89, 5).Value = 153
$ws.Cells.Item(89, 6).Value = 33

$ws.Cells.Item(90, 1).Value = 45222.602939652774
$ws.Cells.Item(90, 2).Value = "Surya Pratama"
$ws.Cells.Item(90, 3).Value = 23123017
$ws.Cells.Item(90, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(90, 5).Value = 178
$ws.Cells.Item(90, 6).Value = 48

$ws.Cells.Item(91, 1).Value = 45222.607119988425
$ws.Cells.Item(91, 2).Value = "Muhammad Reza Fadhila"
$ws.Cells.Item(91, 3).Value = 23123041
$ws.Cells.Item(91, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(91, 5).Value = 173
$ws.Cells.Item(91, 6).Value = 56

$ws.Cells.Item(92, 1).Value = 45222.627195231486
$ws.Cells.Item(92, 2).Value = "Azkaa Satria"
$ws.Cells.Item(92, 3).Value = 23123040
$ws.Cells.Item(92, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(92, 5).Value = 171
$ws.Cells.Item(92, 6).Value = 56

$ws.Cells.Item(93, 1).Value = 45222.68596700231
$ws.Cells.Item(93, 2).Value = "Harits Satriaksa"
$ws.Cells.Item(93, 3).Value = 23123035
$ws.Cells.Item(93, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(93, 5).Value = 174
$ws.Cells.Item(93, 6).Value = 273

$ws.Cells.Item(94, 1).Value = 45222.68635070602
$ws.Cells.Item(94, 2).Value = "Henggar Agung Wirawan"
$ws.Cells.Item(94, 3).Value = 23123024
$ws.Cells.Item(94, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(94, 5).Value = 170
$ws.Cells.Item(94, 6).Value = 61

$ws.Cells.Item(95, 1).Value = 45222.68732729167
$ws.Cells.Item(95, 2).Value = "M Abdurachman Alfatih"
$ws.Cells.Item(95, 3).Value = 23123023
$ws.Cells.Item(95, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(95, 5).Value = 168
$ws.Cells.Item(95, 6).Value = "33 WPM (141CPM)"

$ws.Cells.Item(96, 1).Value = 45222.689265567125
$ws.Cells.Item(96, 2).Value = "Gilang Cahyo Nugroho"
$ws.Cells.Item(96, 3).Value = 23123030
$ws.Cells.Item(96, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(96, 5).Value = 160
$ws.Cells.Item(96, 6).Value = 44

$ws.Cells.Item(97, 1).Value = 45222.69018539352
$ws.Cells.Item(97, 2).Value = "Mokhamad Irfan"
$ws.Cells.Item(97, 3).Value = 23123037
$ws.Cells.Item(97, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(97, 5).Value = "`"185`""
$ws.Cells.Item(97, 6).Value = 56

$ws.Cells.Item(98, 1).Value = 45222.694287962964
$ws.Cells.Item(98, 2).Value = "Wisnu Sri Nugroho"
$ws.Cells.Item(98, 3).Value = 23123014
$ws.Cells.Item(98, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(98, 5).Value = 169
$ws.Cells.Item(98, 6).Value = 42

$ws.Cells.Item(99, 1).Value = 45222.698357905094
$ws.Cells.Item(99, 2).Value = "Dwaldes Bernad"
$ws.Cells.Item(99, 3).Value = 23123029
$ws.Cells.Item(99, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(99, 5).Value = 168
$ws.Cells.Item(99, 6).Value = 50

$ws.Cells.Item(100, 1).Value = 45222.708710439816
$ws.Cells.Item(100, 2).Value = "CHITRA PHRISTIAWAN AJI BHUWANA"
$ws.Cells.Item(100, 3).Value = 23123033
$ws.Cells.Item(100, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(100, 5).Value = 173
$ws.Cells.Item(100, 6).Value = 42

$ws.Cells.Item(101, 1).Value = 45222.72075579861
$ws.Cells.Item(101, 2).Value = "Adhesty Darmayanti Ratulasmar"
$ws.Cells.Item(101, 3).Value = 23123018
$ws.Cells.Item(101, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(101, 5).Value = 155
$ws.Cells.Item(101, 6).Value = 53

$ws.Cells.Item(102, 1).Value = 45222.74390423611
$ws.Cells.Item(102, 2).Value = "Syamsurya Catur Aprian"
$ws.Cells.Item(102, 3).Value = 23123019
$ws.Cells.Item(102, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(102, 5).Value = 176
$ws.Cells.Item(102, 6).Value = 36

$ws.Cells.Item(103, 1).Value = 45222.748589780094
$ws.Cells.Item(103, 2).Value = "Dian Priyatno"
$ws.Cells.Item(103, 3).Value = 23123036
$ws.Cells.Item(103, 4).Value = "Non-regular (PLN)"
$ws.Cells.Item(103, 5).Value = 171
$ws.Cells.Item(103, 6).Value = 32
